{"js": "const body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"text\");\nawait context.sync();\n\n// Locate the two target paragraphs by their text content.\nlet introPara = null;\nlet videoPara = null;\nfor (let i = 0; i < paras.items.length; i++) {\n  const t = paras.items[i].text;\n  if (t === \"THE SMARTCARD INTRODUCTION\") {\n    introPara = paras.items[i];\n  } else if (t.indexOf(\"We\\u2019ve prepared an short video\") === 0) {\n    videoPara = paras.items[i];\n  }\n}\n\nif (!introPara || !videoPara) {\n  throw new Error(\"Could not locate target paragraphs\");\n}\n\nconst startRange = introPara.getRange(\"Start\");\nconst endRange = videoPara.getRange(\"End\");\nconst targetRange = startRange.expandTo(endRange);\n\nconst newBodyXml = `    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading2\"/>\n        <w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/>\n        <w:spacing w:before=\"0\" w:beforeAutospacing=\"0\" w:after=\"150\" w:afterAutospacing=\"0\" w:line=\"264\" w:lineRule=\"atLeast\"/>\n        <w:jc w:val=\"center\"/>\n        <w:textAlignment w:val=\"baseline\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"exo\" w:hAnsi=\"exo\" w:cs=\"open sans\"/>\n          <w:caps/>\n          <w:color w:val=\"343434\"/>\n          <w:spacing w:val=\"15\"/>\n          <w:sz w:val=\"42\"/>\n          <w:szCs w:val=\"42\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"exo\" w:hAnsi=\"exo\" w:cs=\"open sans\"/>\n          <w:caps/>\n          <w:color w:val=\"343434\"/>\n          <w:spacing w:val=\"15\"/>\n          <w:sz w:val=\"42\"/>\n          <w:szCs w:val=\"42\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\">THE SMARTCARD INTRODUCTION</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading2\"/>\n        <w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/>\n        <w:spacing w:before=\"0\" w:beforeAutospacing=\"0\" w:after=\"150\" w:afterAutospacing=\"0\" w:line=\"264\" w:lineRule=\"atLeast\"/>\n        <w:jc w:val=\"center\"/>\n        <w:textAlignment w:val=\"baseline\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"exo\" w:hAnsi=\"exo\"/>\n          <w:caps/>\n          <w:color w:val=\"343434\"/>\n          <w:spacing w:val=\"15\"/>\n          <w:sz w:val=\"42\"/>\n          <w:szCs w:val=\"42\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"exo\" w:hAnsi=\"exo\"/>\n          <w:caps/>\n          <w:color w:val=\"343434\"/>\n          <w:spacing w:val=\"15\"/>\n          <w:sz w:val=\"42\"/>\n          <w:szCs w:val=\"42\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\">SMARTCASH FOR BUSINESS</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading2\"/>\n        <w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/>\n        <w:spacing w:before=\"0\" w:beforeAutospacing=\"0\" w:after=\"150\" w:afterAutospacing=\"0\" w:line=\"264\" w:lineRule=\"atLeast\"/>\n        <w:jc w:val=\"center\"/>\n        <w:textAlignment w:val=\"baseline\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"exo\" w:hAnsi=\"exo\" w:cs=\"open sans\"/>\n          <w:caps/>\n          <w:color w:val=\"343434\"/>\n          <w:spacing w:val=\"15\"/>\n          <w:sz w:val=\"42\"/>\n          <w:szCs w:val=\"42\"/>\n        </w:rPr>\n      </w:pPr>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading1\"/>\n        <w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/>\n        <w:spacing w:before=\"0\" w:beforeAutospacing=\"0\" w:after=\"210\" w:afterAutospacing=\"0\" w:line=\"264\" w:lineRule=\"atLeast\"/>\n        <w:jc w:val=\"center\"/>\n        <w:textAlignment w:val=\"baseline\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"exo\" w:hAnsi=\"exo\"/>\n          <w:caps/>\n          <w:color w:val=\"343434\"/>\n          <w:spacing w:val=\"15\"/>\n          <w:sz w:val=\"51\"/>\n          <w:szCs w:val=\"51\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"exo\" w:hAnsi=\"exo\"/>\n          <w:caps/>\n          <w:color w:val=\"343434\"/>\n          <w:spacing w:val=\"15\"/>\n          <w:sz w:val=\"51\"/>\n          <w:szCs w:val=\"51\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\">SMARTCARD TUTORIAL VIDEOS</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading1\"/>\n        <w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/>\n        <w:spacing w:before=\"0\" w:beforeAutospacing=\"0\" w:after=\"210\" w:afterAutospacing=\"0\" w:line=\"264\" w:lineRule=\"atLeast\"/>\n        <w:textAlignment w:val=\"baseline\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"exo\" w:hAnsi=\"exo\"/>\n          <w:caps/>\n          <w:color w:val=\"343434\"/>\n          <w:spacing w:val=\"15\"/>\n          <w:sz w:val=\"51\"/>\n          <w:szCs w:val=\"51\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"exo\" w:hAnsi=\"exo\"/>\n          <w:caps/>\n          <w:color w:val=\"343434\"/>\n          <w:spacing w:val=\"15\"/>\n          <w:sz w:val=\"51\"/>\n          <w:szCs w:val=\"51\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\">DOWNLOAD SMARTPAY APP</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"NormalWeb\"/>\n        <w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/>\n        <w:spacing w:before=\"204\" w:beforeAutospacing=\"0\" w:after=\"204\" w:afterAutospacing=\"0\"/>\n        <w:textAlignment w:val=\"baseline\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"open sans\" w:hAnsi=\"open sans\" w:cs=\"open sans\"/>\n          <w:color w:val=\"252525\"/>\n          <w:sz w:val=\"21\"/>\n          <w:szCs w:val=\"21\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"open sans\" w:hAnsi=\"open sans\" w:cs=\"open sans\"/>\n          <w:color w:val=\"252525\"/>\n          <w:sz w:val=\"21\"/>\n          <w:szCs w:val=\"21\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\">Accept SmartCash as a payment option in your business with zero fees using a simple SmartPay app. The SmartPay app is available for use anywhere in the world, all it requires is access to the internet.</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading2\"/>\n        <w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/>\n        <w:spacing w:before=\"0\" w:beforeAutospacing=\"0\" w:after=\"150\" w:afterAutospacing=\"0\" w:line=\"264\" w:lineRule=\"atLeast\"/>\n        <w:textAlignment w:val=\"baseline\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"exo\" w:hAnsi=\"exo\"/>\n          <w:caps/>\n          <w:color w:val=\"343434\"/>\n          <w:spacing w:val=\"15\"/>\n          <w:sz w:val=\"42\"/>\n          <w:szCs w:val=\"42\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"exo\" w:hAnsi=\"exo\"/>\n          <w:caps/>\n          <w:color w:val=\"343434\"/>\n          <w:spacing w:val=\"15\"/>\n          <w:sz w:val=\"42\"/>\n          <w:szCs w:val=\"42\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\">SMARTCARD SHOP</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"NormalWeb\"/>\n        <w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/>\n        <w:spacing w:before=\"204\" w:beforeAutospacing=\"0\" w:after=\"204\" w:afterAutospacing=\"0\"/>\n        <w:textAlignment w:val=\"baseline\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"open sans\" w:hAnsi=\"open sans\" w:cs=\"open sans\"/>\n          <w:color w:val=\"252525\"/>\n          <w:sz w:val=\"21\"/>\n          <w:szCs w:val=\"21\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"open sans\" w:hAnsi=\"open sans\" w:cs=\"open sans\"/>\n          <w:color w:val=\"252525\"/>\n          <w:sz w:val=\"21\"/>\n          <w:szCs w:val=\"21\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\">An online store where you can purchase physical SmartCards.</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading2\"/>\n        <w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/>\n        <w:spacing w:before=\"0\" w:beforeAutospacing=\"0\" w:after=\"150\" w:afterAutospacing=\"0\" w:line=\"264\" w:lineRule=\"atLeast\"/>\n        <w:textAlignment w:val=\"baseline\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"exo\" w:hAnsi=\"exo\"/>\n          <w:caps/>\n          <w:color w:val=\"343434\"/>\n          <w:spacing w:val=\"15\"/>\n          <w:sz w:val=\"42\"/>\n          <w:szCs w:val=\"42\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"exo\" w:hAnsi=\"exo\"/>\n          <w:caps/>\n          <w:color w:val=\"343434\"/>\n          <w:spacing w:val=\"15\"/>\n          <w:sz w:val=\"42\"/>\n          <w:szCs w:val=\"42\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\">READY TO GET STARTED?</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading2\"/>\n        <w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/>\n        <w:spacing w:before=\"0\" w:beforeAutospacing=\"0\" w:after=\"150\" w:afterAutospacing=\"0\" w:line=\"264\" w:lineRule=\"atLeast\"/>\n        <w:jc w:val=\"center\"/>\n        <w:textAlignment w:val=\"baseline\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"exo\" w:hAnsi=\"exo\" w:cs=\"open sans\"/>\n          <w:caps/>\n          <w:color w:val=\"343434\"/>\n          <w:spacing w:val=\"15\"/>\n          <w:sz w:val=\"42\"/>\n          <w:szCs w:val=\"42\"/>\n        </w:rPr>\n      </w:pPr>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/>\n        <w:spacing w:line=\"396\" w:lineRule=\"atLeast\"/>\n        <w:textAlignment w:val=\"baseline\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"inherit\" w:hAnsi=\"inherit\" w:cs=\"open sans\"/>\n          <w:color w:val=\"252525\"/>\n          <w:sz w:val=\"21\"/>\n          <w:szCs w:val=\"21\"/>\n        </w:rPr>\n      </w:pPr>\n    </w:p>\n`;\n\nconst ooxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' + newBodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>';\n\ntargetRange.insertOoxml(ooxml, \"Replace\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$introIndex = -1\n$videoIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t.StartsWith(\"THE SMARTCARD INTRODUCTION\")) {\n        $introIndex = $i\n    } elseif ($t.StartsWith(\"We\u2019ve prepared an short video\")) {\n        $videoIndex = $i\n    }\n}\n\nif ($introIndex -eq -1 -or $videoIndex -eq -1) {\n    throw \"Could not locate target paragraphs\"\n}\n\n$startPara = $d.Paragraphs.Item($introIndex)\n$endPara = $d.Paragraphs.Item($videoIndex)\n$targetRange = $d.Range($startPara.Range.Start, $endPara.Range.End)\n\n$newBodyXml = '    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading2\"/>\n        <w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/>\n        <w:spacing w:before=\"0\" w:beforeAutospacing=\"0\" w:after=\"150\" w:afterAutospacing=\"0\" w:line=\"264\" w:lineRule=\"atLeast\"/>\n        <w:jc w:val=\"center\"/>\n        <w:textAlignment w:val=\"baseline\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"exo\" w:hAnsi=\"exo\" w:cs=\"open sans\"/>\n          <w:caps/>\n          <w:color w:val=\"343434\"/>\n          <w:spacing w:val=\"15\"/>\n          <w:sz w:val=\"42\"/>\n          <w:szCs w:val=\"42\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"exo\" w:hAnsi=\"exo\" w:cs=\"open sans\"/>\n          <w:caps/>\n          <w:color w:val=\"343434\"/>\n          <w:spacing w:val=\"15\"/>\n          <w:sz w:val=\"42\"/>\n          <w:szCs w:val=\"42\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\">THE SMARTCARD INTRODUCTION</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading2\"/>\n        <w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/>\n        <w:spacing w:before=\"0\" w:beforeAutospacing=\"0\" w:after=\"150\" w:afterAutospacing=\"0\" w:line=\"264\" w:lineRule=\"atLeast\"/>\n        <w:jc w:val=\"center\"/>\n        <w:textAlignment w:val=\"baseline\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"exo\" w:hAnsi=\"exo\"/>\n          <w:caps/>\n          <w:color w:val=\"343434\"/>\n          <w:spacing w:val=\"15\"/>\n          <w:sz w:val=\"42\"/>\n          <w:szCs w:val=\"42\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"exo\" w:hAnsi=\"exo\"/>\n          <w:caps/>\n          <w:color w:val=\"343434\"/>\n          <w:spacing w:val=\"15\"/>\n          <w:sz w:val=\"42\"/>\n          <w:szCs w:val=\"42\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\">SMARTCASH FOR BUSINESS</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading2\"/>\n        <w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/>\n        <w:spacing w:before=\"0\" w:beforeAutospacing=\"0\" w:after=\"150\" w:afterAutospacing=\"0\" w:line=\"264\" w:lineRule=\"atLeast\"/>\n        <w:jc w:val=\"center\"/>\n        <w:textAlignment w:val=\"baseline\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"exo\" w:hAnsi=\"exo\" w:cs=\"open sans\"/>\n          <w:caps/>\n          <w:color w:val=\"343434\"/>\n          <w:spacing w:val=\"15\"/>\n          <w:sz w:val=\"42\"/>\n          <w:szCs w:val=\"42\"/>\n        </w:rPr>\n      </w:pPr>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading1\"/>\n        <w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/>\n        <w:spacing w:before=\"0\" w:beforeAutospacing=\"0\" w:after=\"210\" w:afterAutospacing=\"0\" w:line=\"264\" w:lineRule=\"atLeast\"/>\n        <w:jc w:val=\"center\"/>\n        <w:textAlignment w:val=\"baseline\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"exo\" w:hAnsi=\"exo\"/>\n          <w:caps/>\n          <w:color w:val=\"343434\"/>\n          <w:spacing w:val=\"15\"/>\n          <w:sz w:val=\"51\"/>\n          <w:szCs w:val=\"51\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"exo\" w:hAnsi=\"exo\"/>\n          <w:caps/>\n          <w:color w:val=\"343434\"/>\n          <w:spacing w:val=\"15\"/>\n          <w:sz w:val=\"51\"/>\n          <w:szCs w:val=\"51\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\">SMARTCARD TUTORIAL VIDEOS</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading1\"/>\n        <w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/>\n        <w:spacing w:before=\"0\" w:beforeAutospacing=\"0\" w:after=\"210\" w:afterAutospacing=\"0\" w:line=\"264\" w:lineRule=\"atLeast\"/>\n        <w:textAlignment w:val=\"baseline\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"exo\" w:hAnsi=\"exo\"/>\n          <w:caps/>\n          <w:color w:val=\"343434\"/>\n          <w:spacing w:val=\"15\"/>\n          <w:sz w:val=\"51\"/>\n          <w:szCs w:val=\"51\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"exo\" w:hAnsi=\"exo\"/>\n          <w:caps/>\n          <w:color w:val=\"343434\"/>\n          <w:spacing w:val=\"15\"/>\n          <w:sz w:val=\"51\"/>\n          <w:szCs w:val=\"51\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\">DOWNLOAD SMARTPAY APP</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"NormalWeb\"/>\n        <w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/>\n        <w:spacing w:before=\"204\" w:beforeAutospacing=\"0\" w:after=\"204\" w:afterAutospacing=\"0\"/>\n        <w:textAlignment w:val=\"baseline\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"open sans\" w:hAnsi=\"open sans\" w:cs=\"open sans\"/>\n          <w:color w:val=\"252525\"/>\n          <w:sz w:val=\"21\"/>\n          <w:szCs w:val=\"21\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"open sans\" w:hAnsi=\"open sans\" w:cs=\"open sans\"/>\n          <w:color w:val=\"252525\"/>\n          <w:sz w:val=\"21\"/>\n          <w:szCs w:val=\"21\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\">Accept SmartCash as a payment option in your business with zero fees using a simple SmartPay app. The SmartPay app is available for use anywhere in the world, all it requires is access to the internet.</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading2\"/>\n        <w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/>\n        <w:spacing w:before=\"0\" w:beforeAutospacing=\"0\" w:after=\"150\" w:afterAutospacing=\"0\" w:line=\"264\" w:lineRule=\"atLeast\"/>\n        <w:textAlignment w:val=\"baseline\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"exo\" w:hAnsi=\"exo\"/>\n          <w:caps/>\n          <w:color w:val=\"343434\"/>\n          <w:spacing w:val=\"15\"/>\n          <w:sz w:val=\"42\"/>\n          <w:szCs w:val=\"42\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"exo\" w:hAnsi=\"exo\"/>\n          <w:caps/>\n          <w:color w:val=\"343434\"/>\n          <w:spacing w:val=\"15\"/>\n          <w:sz w:val=\"42\"/>\n          <w:szCs w:val=\"42\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\">SMARTCARD SHOP</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"NormalWeb\"/>\n        <w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/>\n        <w:spacing w:before=\"204\" w:beforeAutospacing=\"0\" w:after=\"204\" w:afterAutospacing=\"0\"/>\n        <w:textAlignment w:val=\"baseline\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"open sans\" w:hAnsi=\"open sans\" w:cs=\"open sans\"/>\n          <w:color w:val=\"252525\"/>\n          <w:sz w:val=\"21\"/>\n          <w:szCs w:val=\"21\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"open sans\" w:hAnsi=\"open sans\" w:cs=\"open sans\"/>\n          <w:color w:val=\"252525\"/>\n          <w:sz w:val=\"21\"/>\n          <w:szCs w:val=\"21\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\">An online store where you can purchase physical SmartCards.</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading2\"/>\n        <w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/>\n        <w:spacing w:before=\"0\" w:beforeAutospacing=\"0\" w:after=\"150\" w:afterAutospacing=\"0\" w:line=\"264\" w:lineRule=\"atLeast\"/>\n        <w:textAlignment w:val=\"baseline\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"exo\" w:hAnsi=\"exo\"/>\n          <w:caps/>\n          <w:color w:val=\"343434\"/>\n          <w:spacing w:val=\"15\"/>\n          <w:sz w:val=\"42\"/>\n          <w:szCs w:val=\"42\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"exo\" w:hAnsi=\"exo\"/>\n          <w:caps/>\n          <w:color w:val=\"343434\"/>\n          <w:spacing w:val=\"15\"/>\n          <w:sz w:val=\"42\"/>\n          <w:szCs w:val=\"42\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\">READY TO GET STARTED?</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading2\"/>\n        <w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/>\n        <w:spacing w:before=\"0\" w:beforeAutospacing=\"0\" w:after=\"150\" w:afterAutospacing=\"0\" w:line=\"264\" w:lineRule=\"atLeast\"/>\n        <w:jc w:val=\"center\"/>\n        <w:textAlignment w:val=\"baseline\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"exo\" w:hAnsi=\"exo\" w:cs=\"open sans\"/>\n          <w:caps/>\n          <w:color w:val=\"343434\"/>\n          <w:spacing w:val=\"15\"/>\n          <w:sz w:val=\"42\"/>\n          <w:szCs w:val=\"42\"/>\n        </w:rPr>\n      </w:pPr>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/>\n        <w:spacing w:line=\"396\" w:lineRule=\"atLeast\"/>\n        <w:textAlignment w:val=\"baseline\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"inherit\" w:hAnsi=\"inherit\" w:cs=\"open sans\"/>\n          <w:color w:val=\"252525\"/>\n          <w:sz w:val=\"21\"/>\n          <w:szCs w:val=\"21\"/>\n        </w:rPr>\n      </w:pPr>\n    </w:p>\n'\n\n$xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' + '<w:body>' + $newBodyXml + '</w:body>' + '</w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$targetRange.InsertXML($xml)\n"}
